$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.346.42'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '1.865.77'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.69'
$ws.Range('E5').Value = '  +3.42%  '
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.90'
$ws.Range('E8').Value = '  +7.56%  '
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('E10').Value = '  +1.63%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').Value = '2.136.11'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.58'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.682'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.844.77'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').Value = '35.324.56'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.29'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '0.0₃0800'
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.42'
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.27'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.72'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.92'
$ws.Range('E26').Value = '  +26.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.22'
$ws.Range('E27').Value = '  +5.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.79'
$ws.Range('E28').Value = '  +2.01%  '
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('E32').Value = '  +2.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.85'
$ws.Range('E33').Value = '  +27.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.06'
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  +9.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.817'
$ws.Range('E36').Value = '  +17.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.31'
$ws.Range('E37').Value = '  +6.68%  '
$ws.Range('E38').Value = '  +4.19%  '
$ws.Range('E39').Value = '  +4.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '91.00'
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('D41').Value = '1.348.78'
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.24'
$ws.Range('E42').Value = '  +2.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0603'
$ws.Range('E43').Value = '  +15.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.36'
$ws.Range('E44').Value = '  +2.30%  '
$ws.Range('B45').Value = 'Gas'
$ws.Range('C45').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.74'
$ws.Range('E45').Value = '  +50.97%  '
$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.64'
$ws.Range('E47').Value = '  +5.70%  '
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = '2.050.29'
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0687'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.42'
$ws.Range('E51').Value = '  +0.28%  '

Write-Output "Applied 92 cell updates"
